$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{}

$data[2] = @("ECs", "Jag1", "Notch3", "ECs", 3, 1, 10.27784, 30.83352, 0.230301226653591, 0.230301226653591, 2, 0.6666666666666666, 3.684362666666666, 11.053088, 0.0397011572965827, 0.03970115729658269, 37.86728998997333, 340.80560990976, 0.00914322522497016, 0.009143225224970157)
$data[3] = @("ECs", "Jag1", "Notch3", "FAPs", 3, 1, 10.27784, 30.83352, 0.230301226653591, 0.230301226653591, 3, 1, 5.038243666666667, 15.114731, 0.05429001496473517, 0.05429001496473517, 51.78226228701334, 466.0403605831201, 0.01250305704142032, 0.01250305704142032)
$data[4] = @("ECs", "Jag1", "Notch3", "M2", 3, 1, 10.27784, 30.83352, 0.230301226653591, 0.230301226653591, 2, 0.6666666666666666, 0.09716666666666667, 0.2915, 0.001047027523164011, 0.001047027523164011, 0.9986634533333333, 8.987971079999999, 0.000241131722924743, 0.0002411317229247429)
$data[5] = @("ECs", "Jag1", "Notch3", "sCs", 3, 1, 10.27784, 30.83352, 0.230301226653591, 0.230301226653591, 3, 1, 83.98262666666666, 251.94788, 0.9049618002155182, 0.9049618002155182, 863.1599996597332, 7768.4399969376, 0.2084138126642758, 0.2084138126642758)
$data[6] = @("FAPs", "Jag1", "Notch3", "ECs", 3, 1, 14.56812733333333, 43.704382, 0.3264360600001921, 0.326436060000192, 2, 0.6666666666666666, 3.684362666666666, 11.053088, 0.0397011572965827, 0.03970115729658269, 53.67426447017954, 483.0683802316159, 0.01295988936534433, 0.01295988936534433)
$data[7] = @("FAPs", "Jag1", "Notch3", "FAPs", 3, 1, 14.56812733333333, 43.704382, 0.3264360600001921, 0.326436060000192, 3, 1, 5.038243666666667, 15.114731, 0.05429001496473517, 0.05429001496473517, 73.39777527236022, 660.579977451242, 0.01772221858243962, 0.01772221858243962)
$data[8] = @("FAPs", "Jag1", "Notch3", "M2", 3, 1, 14.56812733333333, 43.704382, 0.3264360600001921, 0.326436060000192, 2, 0.6666666666666666, 0.09716666666666667, 0.2915, 0.001047027523164011, 0.001047027523164011, 1.415536372555555, 12.739827353, 0.0003417875393734197, 0.0003417875393734197)
$data[9] = @("FAPs", "Jag1", "Notch3", "sCs", 3, 1, 14.56812733333333, 43.704382, 0.3264360600001921, 0.326436060000192, 3, 1, 83.98262666666666, 251.94788, 0.9049618002155182, 0.9049618002155182, 1223.469599067795, 11011.22639161016, 0.2954121645130348, 0.2954121645130347)
$data[10] = @("M2", "Jag1", "Notch3", "ECs", 3, 1, 2.650137, 7.950411, 0.05938308067649115, 0.05938308067649114, 2, 0.6666666666666666, 3.684362666666666, 11.053088, 0.0397011572965827, 0.03970115729658269, 9.764065824351999, 87.87659241916799, 0.002357577026693035, 0.002357577026693035)
$data[11] = @("M2", "Jag1", "Notch3", "FAPs", 3, 1, 2.650137, 7.950411, 0.05938308067649115, 0.05938308067649114, 3, 1, 5.038243666666667, 15.114731, 0.05429001496473517, 0.05429001496473517, 13.352035956049, 120.168323604441, 0.00322390833857878, 0.00322390833857878)
$data[12] = @("M2", "Jag1", "Notch3", "M2", 3, 1, 2.650137, 7.950411, 0.05938308067649115, 0.05938308067649114, 2, 0.6666666666666666, 0.09716666666666667, 0.2915, 0.001047027523164011, 0.001047027523164011, 0.2575049785, 2.3175448065, 0.00006217571987855518, 0.00006217571987855518)
$data[13] = @("M2", "Jag1", "Notch3", "sCs", 3, 1, 2.650137, 7.950411, 0.05938308067649115, 0.05938308067649114, 3, 1, 83.98262666666666, 251.94788, 0.9049618002155182, 0.9049618002155182, 222.56546628652, 2003.08919657868, 0.05373941959134078, 0.05373941959134077)
$data[14] = @("sCs", "Jag1", "Notch3", "ECs", 3, 1, 17.13170833333333, 51.395125, 0.3838796326697257, 0.3838796326697257, 2, 0.6666666666666666, 3.684362666666666, 11.053088, 0.0397011572965827, 0.03970115729658269, 63.11942659955555, 568.0748393959999, 0.01524046567957517, 0.01524046567957516)
$data[15] = @("sCs", "Jag1", "Notch3", "FAPs", 3, 1, 17.13170833333333, 51.395125, 0.3838796326697257, 0.3838796326697257, 3, 1, 5.038243666666667, 15.114731, 0.05429001496473517, 0.05429001496473517, 86.31372100959723, 776.8234890863752, 0.02084083100229645, 0.02084083100229645)
$data[16] = @("sCs", "Jag1", "Notch3", "M2", 3, 1, 17.13170833333333, 51.395125, 0.3838796326697257, 0.3838796326697257, 2, 0.6666666666666666, 0.09716666666666667, 0.2915, 0.001047027523164011, 0.001047027523164011, 1.664630993055555, 14.9816789375, 0.0004019325409872934, 0.0004019325409872933)
$data[17] = @("sCs", "Jag1", "Notch3", "sCs", 3, 1, 17.13170833333333, 51.395125, 0.3838796326697257, 0.3838796326697257, 3, 1, 83.98262666666666, 251.94788, 0.9049618002155182, 0.9049618002155182, 1438.765865120555, 12948.892786085, 0.3473964034468668, 0.3473964034468668)

foreach ($r in $data.Keys) {
  $rowVals = $data[$r]
  for ($c = 0; $c -lt $rowVals.Length; $c++) {
    $ws.Cells.Item($r, $c + 1).Value = $rowVals[$c]
  }
}